$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on price cells whose new values would otherwise
# be auto-coerced to numeric by Excel (they must remain text, matching the
# original inlineStr cell type).
foreach ($addr in @("D5", "D8", "D9", "D10", "D17", "D20", "D22", "D24", "D33", "D39", "D41", "D46", "D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.956.95"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "1.640.25"
$ws.Range("E4").Value = "  +0.99%  "
$ws.Range("D5").Value = "214.86"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E7").Value = "  +1.02%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "0.0639"
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.255"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "19.68"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "1.867.08"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").Value = "1.617.55"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "0.0₃0762"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "62.63"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "25.966.25"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "194.21"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("D22").Value = "9.95"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("D24").Value = "144.14"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").Value = "3.24"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("E34").Value = "  -2.53%  "
$ws.Range("E35").Value = "  +2.04%  "
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").Value = "1.140.74"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").Value = "2.46"
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").Value = "99.43"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("D44").Value = "1.776.27"
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("E45").Value = "  +8.19%  "
$ws.Range("D46").Value = "56.57"
$ws.Range("E46").Value = "  +1.54%  "
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "7.63"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("E51").Value = "  -0.57%  "
